# Updated RegisterPage and added Login module
$wb = $excel.ActiveWorkbook
$register = $wb.Worksheets.Item(1)

# --- Register sheet tweaks -------------------------------------------------
# Best-fit style width for column E (mirrors existing bestFit cols on the sheet)
$register.Columns.Item(5).ColumnWidth = 16.109375

# --- Add the new Login sheet, placed right after Register ------------------
$login = $wb.Worksheets.Add($null, $register)
$login.Name = "Login"

# Column widths
$login.Columns.Item(1).ColumnWidth = 21.5546875
$login.Columns.Item(2).ColumnWidth = 14.88671875
$login.Columns.Item(3).ColumnWidth = 15.33203125
$login.Columns.Item(4).ColumnWidth = 14.44140625

# --- Cell values, written in an order that keeps username/password e-mails
# --- (which are re-used across rows) trailing, same as the source data ----
$login.Range("A1").Value = "Username"
$login.Range("B1").Value = "Password"
$login.Range("C1").Value = "Scenario Type"
$login.Range("D1").Value = "Expected Result"

$login.Range("B2").Value = 123456
$login.Range("C2").Value = "Valid"
$login.Range("D2").Value = "Login successful"

$login.Range("A3").Value = "invalid_user@mail.com"
$login.Range("B3").Value = "WrongPass"
$login.Range("C3").Value = "Invalid credentials"
$login.Range("D3").Value = "Error: Login was unsuccessful"

$login.Range("B4").Value = 123456
$login.Range("C4").Value = "Blank username"
$login.Range("D4").Value = "Error: Please enter email"

$login.Range("C5").Value = "Blank password"
$login.Range("D5").Value = "Error: Please enter password"

$login.Range("A6").Value = "Invalid credentials"
$login.Range("B6").Value = 123456
$login.Range("C6").Value = "Invalid email format"
$login.Range("D6").Value = "Error: Wrong email format"

$login.Range("C7").Value = "Both fields blank"
$login.Range("D7").Value = "Error: Please enter credentials"

$login.Range("A2").Value = "email1@gamil.com"
$login.Range("A5").Value = "email1@gamil.com"

# --- Header formatting: bold, centered, wrapped -----------------------------
$login.Range("A1").Font.Bold = $true
$login.Range("A1").HorizontalAlignment = -4108
$login.Range("A1").VerticalAlignment = -4108
$login.Range("A1").WrapText = $true
$login.Range("A1").Copy() | Out-Null
$login.Range("B1:D1").PasteSpecial(-4122) | Out-Null

# --- Body formatting: vertically centered, wrapped --------------------------
$login.Range("B2").VerticalAlignment = -4108
$login.Range("B2").WrapText = $true
$login.Range("B2").Copy() | Out-Null
$login.Range("C2:D2").PasteSpecial(-4122) | Out-Null
$login.Range("A3:D4").PasteSpecial(-4122) | Out-Null
$login.Range("B5:D5").PasteSpecial(-4122) | Out-Null
$login.Range("A6:D7").PasteSpecial(-4122) | Out-Null

# --- Hyperlink cells (Username column on the "Valid" rows) ------------------
$login.Range("A2").Value = "email1@gamil.com"
$login.Range("B2").Copy() | Out-Null
$login.Range("A2").PasteSpecial(-4122) | Out-Null
$login.Hyperlinks.Add($login.Range("A2"), "mailto:email1@gamil.com") | Out-Null

$login.Range("A5").Value = "email1@gamil.com"
$login.Range("A2").Copy() | Out-Null
$login.Range("A5").PasteSpecial(-4122) | Out-Null
$login.Hyperlinks.Add($login.Range("A5"), "mailto:email1@gamil.com") | Out-Null

$excel.CutCopyMode = $false

# --- Row heights -------------------------------------------------------------
$login.Rows.Item(3).RowHeight = 28.8
$login.Rows.Item(4).RowHeight = 28.8
$login.Rows.Item(5).RowHeight = 28.8
$login.Rows.Item(6).RowHeight = 28.8
$login.Rows.Item(7).RowHeight = 43.2

# --- Activate Login sheet & set selection ------------------------------------
$login.Activate()
$login.Range("F5").Select() | Out-Null
